# Daily attendance processing - 2025-12-09 05:28:08
#
# Normalizes the "Recorded By" (column G) list on the "Session Analysis
# Results" sheet so that the automated "System" recorder is always listed
# first among the comma-separated contributors for a session, instead of
# trailing after the user/email entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value()

    if ($val -eq $null) { continue }

    if ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value = "system, System, backup@backdoor.com"
        $changed++
    } elseif ($val -eq "backup@backdoor.com, System") {
        $cell.Value = "System, backup@backdoor.com"
        $changed++
    } elseif ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
        $changed++
    }
}

Write-Host ("Recorded By entries reordered: " + $changed)
